# Offer service test improvements:
#  - valid action code for savings with loan simulations
#  - for methods GetBuildingSavingsDepositSchedule, GetBuildingSavingsPaymentSchedule
#    check if balances fits to response of SimulateBuildingSavings

$wb = $excel.ActiveWorkbook

$wsOffer = $wb.Worksheets.Item("01-OfferService")
$wsPokus = $wb.Worksheets.Item("01-OfferServicePokus")

# --- 01-OfferService: LoanActionCode column (H) goes from 4 to 8 for every
#     data row (rows 2-10) ---
for ($r = 2; $r -le 10; $r++) {
    $wsOffer.Cells.Item($r, 8).Value = 8
}

# --- 01-OfferServicePokus: this sheet only has a single data row, the
#     shared formula that used to span F2:G10 should now only cover F2:G2 ---
$wsPokus.Range("F2:G2").Formula = '="false"'

# --- Window / active-sheet bookkeeping: "01-OfferService" becomes the
#     active (selected) sheet/tab, with H10 as the active selected cell;
#     "01-OfferServicePokus" loses tabSelected and keeps its own selection ---
$wsOffer.Activate()
$wsOffer.Range("H10").Select()

# Best-effort: reflect the new window size recorded in the workbook view.
try {
    $win = $excel.ActiveWindow
    $win.Width = 45360
    $win.Height = 23730
} catch {
}

$wb.Save()
